# Commit: Added new test case and also added logic to read multiple sheets
# from excel in excelutils.
#
# - RUNMANAGER (sheet1): new test row (row 7) for the new test case.
# - DATA (sheet2): duplicate rows (10 & 11) for the "change password" test
#   removed now that sheet-reading logic covers USERACCOUNTMANAGEMENTDATA
#   separately.
# - USERACCOUNTMANAGEMENTDATA (sheet3): gains the data rows for the new
#   test case (rows 4 & 5).

$wb = $excel.ActiveWorkbook
$wsRun  = $wb.Worksheets.Item(1)   # RUNMANAGER
$wsData = $wb.Worksheets.Item(2)   # DATA
$wsUser = $wb.Worksheets.Item(3)   # USERACCOUNTMANAGEMENTDATA

# --- 1. RUNMANAGER: add the new test case row (row 7) --------------------
# New shared string "verifyThatUserCannotChangePasswordWithInvalidCurrentPassword"
# is introduced here first, so it lands at the next shared-string index.
$wsRun.Range("A7").Value = "verifyThatUserCannotChangePasswordWithInvalidCurrentPassword"
$wsRun.Range("B7").Value = "To check this test is executed"
$wsRun.Range("C7").Value = "yes"

# --- 2. USERACCOUNTMANAGEMENTDATA: add data rows for the new test case ---
$wsUser.Range("A4").Value = "verifyThatUserCannotChangePasswordWithInvalidCurrentPassword"
$wsUser.Range("B4").Value = "yes"
$wsUser.Range("C4").Value = "Admin"
$wsUser.Range("D4").Value = "admin123"
$wsUser.Range("E4").Value = "Sunil"
$wsUser.Range("F4").Value = "chrome"
$wsUser.Range("G4").Value = "admin112"
$wsUser.Range("H4").Value = "admin111"
$wsUser.Range("I4").Value = "admin111"

$wsUser.Range("A5").Value = "verifyThatUserCannotChangePasswordWithInvalidCurrentPassword"
$wsUser.Range("B5").Value = "yes"
$wsUser.Range("C5").Value = "Admin"
$wsUser.Range("D5").Value = "admin123"
$wsUser.Range("E5").Value = "Sunil"
$wsUser.Range("F5").Value = "firefox"
$wsUser.Range("G5").Value = "admin112"
$wsUser.Range("H5").Value = "admin111"
$wsUser.Range("I5").Value = "admin111"

# Fit column A now that it holds the longer test name.
$wsUser.Columns("A:A").AutoFit()

# --- 3. RUNMANAGER: priority/count for the new row, kept as text (quote- -
#        prefixed) to match the other rows in the priority/count columns.
$wsRun.Range("D7").Value = "'6"
$wsRun.Range("E7").Value = "'1"

# --- 4. DATA: drop the now-duplicated rows 10 & 11 ------------------------
$wsData.Rows("10:11").Delete()

# --- 5. Sheet view / selection bookkeeping --------------------------------
$wsData.Activate()
[void]$wsData.Range("A15").Select()

$wsUser.Activate()
[void]$wsUser.Range("B6").Select()

$wsUser.PageSetup.Orientation = 1

$wsRun.Activate()
[void]$wsRun.Range("D8").Select()
